$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "29.913.75"
Set-TextValue "E2" "  -0.11%  "
Set-TextValue "D3" "1.876.08"
Set-TextValue "E3" "  -0.90%  "
Set-TextValue "D4" "0.9984"
Set-TextValue "E4" "  -0.19%  "
Set-TextValue "D5" "0.7457"
Set-TextValue "E5" "  -4.29%  "
Set-TextValue "D6" "242.56"
Set-TextValue "E6" "  -0.60%  "
Set-TextValue "D7" "0.9990"
Set-TextValue "E7" "  -0.12%  "
Set-TextValue "E8" "  +0.83%  "
Set-TextValue "D9" "0.07207"
Set-TextValue "E9" "  -1.05%  "
Set-TextValue "D10" "24.81"
Set-TextValue "E10" "  -4.10%  "
Set-TextValue "D11" "0.08403"
Set-TextValue "E11" "  -3.33%  "
Set-TextValue "D12" "0.7516"
Set-TextValue "E12" "  -3.02%  "
Set-TextValue "D13" "5.428"
Set-TextValue "E13" "  +0.21%  "
Set-TextValue "D14" "1.858.16"
Set-TextValue "E14" "  -7.56%  "
Set-TextValue "D15" "92.66"
Set-TextValue "E15" "  -2.02%  "
Set-TextValue "D16" "29.903.34"
Set-TextValue "E16" "  -0.16%  "
Set-TextValue "D17" "6.097"
Set-TextValue "E17" "  -1.80%  "
Set-TextValue "D18" "13.58"
Set-TextValue "E18" "  -2.51%  "
Set-TextValue "D19" "243.84"
Set-TextValue "E19" "  -0.92%  "
Set-TextValue "D20" "0.000007822"
Set-TextValue "E20" "  -0.80%  "
Set-TextValue "D21" "0.9983"
Set-TextValue "E21" "  -0.18%  "
Set-TextValue "D22" "2.122.91"
Set-TextValue "E22" "  -9.50%  "
Set-TextValue "D23" "7.988"
Set-TextValue "E23" "  -3.03%  "
Set-TextValue "D24" "0.9985"
Set-TextValue "E24" "  -0.17%  "
Set-TextValue "D25" "0.1558"
Set-TextValue "E25" "  -6.96%  "
Set-TextValue "D26" "9.285"
Set-TextValue "D27" "165.29"
Set-TextValue "E27" "  +1.18%  "
Set-TextValue "D28" "18.63"
Set-TextValue "E28" "  -1.28%  "
Set-TextValue "D29" "2.037"
Set-TextValue "E29" "  -0.73%  "
Set-TextValue "D30" "1.510"
Set-TextValue "E30" "  +5.41%  "
Set-TextValue "D31" "4.597"
Set-TextValue "D32" "1.528"
Set-TextValue "E32" "  -0.91%  "
Set-TextValue "E33" "  +3.27%  "
Set-TextValue "D34" "0.05319"
Set-TextValue "E34" "  -3.13%  "
Set-TextValue "D35" "1.238"
Set-TextValue "E35" "  -0.82%  "
Set-TextValue "D36" "0.7547"
Set-TextValue "D37" "0.9976"
Set-TextValue "E37" "  -0.80%  "
Set-TextValue "D38" "2.700"
Set-TextValue "E38" "  +0.54%  "
Set-TextValue "D39" "0.01961"
Set-TextValue "E39" "  +0.03%  "
Set-TextValue "D40" "2.755"
Set-TextValue "E40" "  -1.34%  "
Set-TextValue "D41" "0.4541"
Set-TextValue "E41" "  +0.50%  "
Set-TextValue "D42" "1.112.38"
Set-TextValue "E42" "  -0.16%  "
Set-TextValue "E43" "  -0.94%  "
Set-TextValue "D44" "72.45"
Set-TextValue "E44" "  -1.81%  "
Set-TextValue "D45" "0.8559"
Set-TextValue "E45" "  +0.31%  "
Set-TextValue "D46" "1.000"
Set-TextValue "E46" "  +0.07%  "
Set-TextValue "D47" "103.51"
Set-TextValue "E47" "  +0.40%  "
Set-TextValue "D48" "3.109"
Set-TextValue "E48" "  +3.62%  "
Set-TextValue "D49" "7.649"
Set-TextValue "E49" "  +0.54%  "
Set-TextValue "D50" "1.842"
Set-TextValue "E50" "  -2.38%  "
Set-TextValue "D51" "2.020.91"
Set-TextValue "E51" "  -9.13%  "
